$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Insert a new "Other" worksheet right after "Meta" (it becomes the
# second tab, ahead of Phases / ALSFRS-R Score / ... / Right-Hand Grip)
# ------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Meta")
$other = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $meta)
$other.Name = "Other"

# Content: a small Name/Value table with an Onset row below it.
$other.Range("A1").Value = "Name"
$other.Range("B1").Value = "Value"
$other.Range("A2").Value = "Onset"

# Row heights matching the rest of the "grip"-style sheets (12.8pt).
$other.Rows.Item(1).RowHeight = 12.8
$other.Rows.Item(2).RowHeight = 12.8

# Page setup / margins matching the other small data sheets.
$other.PageSetup.PaperSize = 9
$other.PageSetup.LeftMargin = 0.7875 * 72
$other.PageSetup.RightMargin = 0.7875 * 72
$other.PageSetup.TopMargin = 1.05277777777778 * 72
$other.PageSetup.BottomMargin = 1.05277777777778 * 72
$other.PageSetup.HeaderMargin = 0.7875 * 72
$other.PageSetup.FooterMargin = 0.7875 * 72
$other.PageSetup.CenterHeader = "&`"Times New Roman,Regular`"&12&Kffffff&A"
$other.PageSetup.CenterFooter = "&`"Times New Roman,Regular`"&12&KffffffPage &P"

# Selection / activation: "Other" becomes the active sheet with B2 selected.
$other.Range("B2").Select()
$other.Activate()

# ------------------------------------------------------------------
# Bump the header row height on the two grip sheets.
# ------------------------------------------------------------------
$leftGrip = $wb.Worksheets.Item("Left-Hand Grip")
$leftGrip.Rows.Item(1).RowHeight = 13.8

$rightGrip = $wb.Worksheets.Item("Right-Hand Grip")
$rightGrip.Rows.Item(1).RowHeight = 13.8
